$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Clients sheet: give column A an explicit width (new <cols>)
# ------------------------------------------------------------------
$clients = $wb.Worksheets.Item("Clients")
$clients.Columns.Item(1).ColumnWidth = 14.333333333333332

# ------------------------------------------------------------------
# 2. New "Organisations" sheet, placed right after "Practitioners"
#    (Copying K5 gives us the right sheetFormatPr/extLst "shape"
#    with no pre-existing <cols>, then we overwrite its contents.)
# ------------------------------------------------------------------
$practitioners = $wb.Worksheets.Item("Practitioners")
$k5 = $wb.Worksheets.Item("K5")
$k5.Copy([System.Reflection.Missing]::Value, $practitioners)
$ws = $wb.Worksheets.Item("K5 (2)")
$ws.Name = "Organisations"

# Trim the copied K5 shape (4 rows x 12 cols) down to 3 rows x 9 cols
$ws.Rows.Item(4).EntireRow.Delete()
$ws.Columns.Item(10).EntireColumn.Delete()
$ws.Columns.Item(10).EntireColumn.Delete()
$ws.Columns.Item(10).EntireColumn.Delete()

# Row 1 (Version header) already reads A1="Version", B1=1 - keep as is.

# Row 2: column headers
$ws.Range("A2").Value = "organisation_path"
$ws.Range("B2").Value = "organisation_key"
$ws.Range("C2").Value = "organisation_name"
$ws.Range("D2").Value = "organisation_legal_name"
$ws.Range("E2").Value = "organisation_abn"
$ws.Range("F2").Value = "organisation_type"
$ws.Range("G2").Value = "organisation_state"
$ws.Range("H2").Value = "organisation_status"
$ws.Range("I2").Value = "organisation_tags"

# Row 3: sample data
$ws.Range("A3").Value = "PHN999:NFP01"
$ws.Range("B3").Value = "NFP01"
$ws.Range("C3").Value = "Test Provider Organisation NFP1"
$ws.Range("D3").ClearContents()
$ws.Range("E3").Value = 42072953425
$ws.Range("F3").Value = 7
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = 1
$ws.Range("I3").ClearContents()

# Column widths
$ws.Columns.Item(1).ColumnWidth = 18.666666666666668
$ws.Columns.Item(2).ColumnWidth = 17.0
$ws.Columns.Item(3).ColumnWidth = 24.833333333333336
$ws.Columns.Item(4).ColumnWidth = 20.833333333333336
$ws.Columns.Item(5).ColumnWidth = 19.666666666666668

# Page margins: left/right 0.75in, top/bottom 1in, header/footer 0.5in
$ws.PageSetup.LeftMargin = 54
$ws.PageSetup.RightMargin = 54
$ws.PageSetup.TopMargin = 72
$ws.PageSetup.BottomMargin = 72
$ws.PageSetup.HeaderMargin = 36
$ws.PageSetup.FooterMargin = 36

# Active selection on the new sheet
$ws.Range("F4").Select() | Out-Null

# ------------------------------------------------------------------
# 3. Window position (best effort - engine may not persist this)
# ------------------------------------------------------------------
$wb.Windows.Item(1).Left = 540
$wb.Windows.Item(1).Top = 2160
